$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.939.24"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.584.51"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E5").Value = "  -2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.46"
$ws.Range("E6").Value = "  -4.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.577.83"
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.626"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.653"
$ws.Range("E11").Value = "  -3.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.47"
$ws.Range("E12").Value = "  -5.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000309"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.60"
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.163.75"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.71"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.602.18"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.004.86"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.57"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.27"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("E24").Value = "  -8.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.72"
$ws.Range("E25").Value = "  +4.06%  "
$ws.Range("E26").Value = "  -3.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.14"
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("E28").Value = "  -7.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.29"
$ws.Range("E29").Value = "  -3.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.97"
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.63"
$ws.Range("E31").Value = "  -3.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "66.83"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.10"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.116"
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "573.04"
$ws.Range("E35").Value = "  -9.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.10"
$ws.Range("E36").Value = "  +7.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.26"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0803"
$ws.Range("E39").Value = "  -4.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.32"
$ws.Range("E40").Value = "  +13.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.394"
$ws.Range("E41").Value = "  -5.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.55"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("E43").Value = "  -7.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.00"
$ws.Range("E44").Value = "  -5.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.230.83"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.44"
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0437"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.68"
$ws.Range("E48").Value = "  +4.41%  "
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  -4.99%  "
